$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 4
$ws.Range("A5").Value = 5

$ws.Range("A6").Select()
